$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.809.83'
$ws.Range('E2').Value = '  +1.31%  '
$ws.Range('D3').Value = '1.650.79'
$ws.Range('E3').Value = '  +1.57%  '
$ws.Range('E4').Value = '  +0.54%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '216.64'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +1.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.507'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.39%  '
$ws.Range('E7').Value = '  +0.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.252'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +2.00%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0628'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.94%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '19.21'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +2.85%  '
$ws.Range('E11').Value = '  +0.30%  '
$ws.Range('D12').Value = '1.878.81'
$ws.Range('E12').Value = '  +1.45%  '
$ws.Range('D13').Value = '1.652.33'
$ws.Range('E13').Value = '  +1.67%  '
$ws.Range('E14').Value = '  +1.77%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.533'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +2.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '65.55'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.07%  '
$ws.Range('D17').Value = '26.812.73'
$ws.Range('E17').Value = '  +1.17%  '
$ws.Range('D18').Value = '0.0₃0745'
$ws.Range('E18').Value = '  +0.96%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '219.05'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +2.56%  '
$ws.Range('E20').Value = '  +0.52%  '
$ws.Range('E21').Value = '  +2.20%  '
$ws.Range('E22').Value = '  +0.75%  '
$ws.Range('E23').Value = '  +15.35%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '9.52'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +2.66%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '146.90'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -1.05%  '
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('E27').Value = '  +0.55%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.13'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +4.44%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.76'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +1.87%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0517'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +2.05%  '
$ws.Range('E31').Value = '  +1.86%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.36'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.38%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.02'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.47%  '
$ws.Range('D34').Value = '1.283.30'
$ws.Range('E34').Value = '  +4.80%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.55'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +3.77%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.44'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +2.56%  '
$ws.Range('E37').Value = '  +3.71%  '
$ws.Range('E38').Value = '  +6.73%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.830'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +4.76%  '
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.816'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +3.07%  '
$ws.Range('E42').Value = '  -0.76%  '
$ws.Range('E43').Value = '  +2.79%  '
$ws.Range('D44').Value = '1.788.71'
$ws.Range('E44').Value = '  +1.57%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '92.10'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -0.61%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '59.93'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +9.60%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.61'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.23%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0516'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +1.26%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.80'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.41%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0971'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +2.35%  '
$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.409'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +0.46%  '
